# Update "想去人数" (want-to-go count) figures in column F
# for the 展览 (Exhibition) and 全部类型 (All Types) sheets.
# These two sheets mirror the same event listing data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 153
    3  = 7072
    4  = 4316
    5  = 67
    10 = 63
    11 = 57
    13 = 616
    14 = 102
    15 = 51
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
